$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.109582662582397
$ws.Range("B1").Value = 1.998434662818909
$ws.Range("D1").Value = 0.9894780516624451
$ws.Range("E1").Value = 1.083810567855835
